# Pick Place for Near-E-Field-Probe - JLCPCB.xlsx
#
# - Replaced connector by the reference: C1509219 - Female header
# - Modified boards to improve manufacturing and asembly: increased connector pad sizes
#
# Concretely (per the recorded OOXML diff) this updates the "Mid X" value
# for part J1 from 60.198mm to 60.1980mm, reflecting the updated pad/
# footprint placement precision after swapping in the new connector.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds designator J1; column B is "Mid X".
$ws.Range("B2").Value2 = "60.1980mm"

# Leave the sheet with the same active selection recorded after the edit.
$ws.Range("F9").Select()
